$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '75.670.82'
$ws.Range('E2').Value = '  +1.99%  '

$ws.Range('D3').Value = '2.828.98'
$ws.Range('E3').Value = '  +6.95%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '191.28'
$ws.Range('E5').Value = '  +2.87%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '595.78'
$ws.Range('E6').Value = '  +2.22%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('E8').Value = '  +3.23%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.191'
$ws.Range('E9').Value = '  +0.32%  '

$ws.Range('D10').Value = '2.829.10'
$ws.Range('E10').Value = '  +7.12%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.383'
$ws.Range('E11').Value = '  +8.35%  '

$ws.Range('E12').Value = '  -1.92%  '

$ws.Range('E13').Value = '  +4.64%  '

$ws.Range('D14').Value = '3.347.65'
$ws.Range('E14').Value = '  +6.21%  '

$ws.Range('D15').Value = '75.525.44'
$ws.Range('E15').Value = '  +1.88%  '

$ws.Range('E16').Value = '  +1.48%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.04'
$ws.Range('E17').Value = '  +3.34%  '

$ws.Range('D18').Value = '2.819.93'
$ws.Range('E18').Value = '  +6.25%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.03'
$ws.Range('E19').Value = '  -2.50%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.34'
$ws.Range('E20').Value = '  +4.21%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '379.67'
$ws.Range('E21').Value = '  +3.24%  '

$ws.Range('E22').Value = '  +2.53%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.11'
$ws.Range('E23').Value = '  +1.77%  '

$ws.Range('E24').Value = '  +0.06%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '71.20'

$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.19'
$ws.Range('E26').Value = '  +2.98%  '

$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.82'
$ws.Range('E27').Value = '  +6.28%  '

$ws.Range('E28').Value = '  +6.38%  '

$ws.Range('E29').Value = '  +12.63%  '

$ws.Range('E30').Value = '  -0.04%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.41'
$ws.Range('E31').Value = '  +2.37%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '516.98'
$ws.Range('E32').Value = '  +0.24%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.71'
$ws.Range('E33').Value = '  +1.52%  '

$ws.Range('E34').Value = '  +4.56%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.19%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '164.84'
$ws.Range('E36').Value = '  +1.30%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '19.90'
$ws.Range('E37').Value = '  +4.24%  '

$ws.Range('E38').Value = '  +0.70%  '

$ws.Range('E39').Value = '  +0.51%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '183.75'
$ws.Range('E40').Value = '  +12.30%  '

$ws.Range('E41').Value = '  -0.01%  '

$ws.Range('E42').Value = '  +5.41%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.04'
$ws.Range('E43').Value = '  +3.24%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.68'
$ws.Range('E44').Value = '  +1.80%  '

$ws.Range('E45').Value = '  +3.39%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '40.00'
$ws.Range('E46').Value = '  +2.70%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0873'
$ws.Range('E47').Value = '  +3.42%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.36'
$ws.Range('E48').Value = '  +1.17%  '

$ws.Range('E49').Value = '  +8.93%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.75'
$ws.Range('E50').Value = '  +4.01%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.647'
$ws.Range('E51').Value = '  +10.15%  '
